$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.166.82'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '1.828.47'
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9986'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6155'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07336'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2902'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.18'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07636'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("D12").Value = '1.841.42'
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.973'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6711'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.41'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008955'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.851'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.88%  '
$ws.Range("D18").Value = '29.163.10'
$ws.Range("E18").Value = '  +0.44%  '
$ws.Range("D19").Value = '2.082.21'
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '235.86'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.44%  '
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.383'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9995'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.78'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.523'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.39%  '
$ws.Range("E27").Value = '  -1.98%  '
$ws.Range("E28").Value = '  -1.26%  '
$ws.Range("E29").Value = '  -0.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05833'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.72%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.076'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.69%  '
$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.219'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.089'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.845'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.133'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.64%  '
$ws.Range("E36").Value = '  -3.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.612'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.71%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.857'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.10%  '
$ws.Range("D39").Value = '1.221.91'
$ws.Range("E39").Value = '  +0.75%  '
$ws.Range("E40").Value = '  -1.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.173'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.97%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8991'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("D44").Value = '2.003.54'
$ws.Range("E44").Value = '  +1.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.77'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.27%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.52'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5045'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.155'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.70%  '
$ws.Range("E49").Value = '  -0.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4024'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1163'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.11%  '
